$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ClientData")

# Update the contact name
$ws.Range("A2").Value = "Erasmus Hartman II"

# Update the VAT number
$ws.Range("L2").Value = "GB21257622945"

# Update the LfrDropDown value
$ws.Range("M2").Value = "NL LMLOG Fiscal Rep (Import)"

# Update the sheet view selection (the engine does not persist topLeftCell,
# only the active selection, so move the selection to L9 as in the edit)
[void]$ws.Range("L9").Select()

# Adjust style of O2 to match the style used by E2/F2/G2 (fontId 4, applyAlignment, no explicit alignment)
$ws.Range("O2").Style = $ws.Range("E2").Style
